$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append after the current last row (170)
$data = @(
    @("2025-07-07", "eaux souterraines", 118, 1),
    @("2025-07-07", "ruissellement",     119, 1),
    @("2025-07-07", "eaux souterraines", 119, 1),
    @("2025-07-07", "eaux de surface",   119, 1)
)

$startRow = 171
$endRow = $startRow + $data.Count - 1

# The "Date" column holds values that look like dates (e.g. "2025-07-07")
# but must stay as plain text, matching the rest of the column. Force the
# range to Text format before assigning so Excel doesn't auto-convert the
# strings into date serials, then drop the formatting again so the new
# cells don't pick up an explicit style (matching the unstyled cells
# already in the sheet).
$dateRange = $ws.Range("A$startRow" + ":A$endRow")
$dateRange.NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$dateRange.ClearFormats()
